$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 value from "*" to "***"
$ws.Range("D2").Value = "***"

# Update selection to D2
$ws.Range("D2").Select()
